$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
  "fällen","arten","treiben","bellen","sichern","sprengen","sterben","lesen",
  "sperren","gründen","schlucken","mauern","graben","bergen","trennen","siegen",
  "flehen","erben","zögern","schwingen","heulen","klettern","enden","hören",
  "hauen","liefern","zeigen","klingen","drehen","stammen","folgen","flüchten",
  "warnen","heben","fallen","wehtun","brauchen","stecken","formen","tropfen",
  "schwören","schmecken","ehren","ärgern","kichern","schrecken","fahren","suchen",
  "dringen","kosten","spüren","wenden","schreiten","rufen","streichen","münzen",
  "wundern","bluten","kümmern","rasen","seufzen","achten","runden","spielen",
  "schulden","backen","tollen","wüten","greifen","pflanzen","biegen","malen",
  "gelten","scheinen","wachsen","bauen","pfeifen","fließen","wirken","saufen",
  "sinken","werfen","trauen","planen","betteln","führen","fangen","lockern",
  "sorgen","reizen","töten","altern","mögen","zielen","jubeln","stehlen",
  "lügen","weichen","quälen","kehren","schwächen","spinnen","dienen","fischen",
  "heilen","helfen","decken","scheitern","jagen","liegen","freuen","platzen",
  "geben","grüßen","machen","ändern","bitten","filmen","feiern","knarren",
  "äußern","zünden","boxen","irren","schenken","loben","räumen","werden"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $words[$i]
}
